$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3: new weekly status entry (22.04.2014 / serial 41751) ---
$ws.Range("A3").Value = 41751
$ws.Range("A3").NumberFormat = "m/d/yy"
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").WrapText = $true

$ws.Range("I3").Value = "1.Війна з ajax, який постійно і нагло викидує Ajax error 2.Розробка CSE Dashboard (його готовність на перерішній час складає 7%)"

# B3:I3 share the bordered + wrap-text style used by the rest of row 2/3
$ws.Range("B3:I3").WrapText = $true

# Row 2 header/status cells (C2:I2) also pick up wrap-text in this edit
$ws.Range("C2:I2").WrapText = $true
$ws.Range("A2").WrapText = $true

$ws.Rows.Item(3).RowHeight = 120

# --- View cosmetics: zoom + active selection moved to I4 after data entry ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("I4").Select()
